# Add a new "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled like the other header cells (bold, bordered,
# centered) by copying G1's format onto H1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell for row 2 (no special formatting, like the other data cells).
$ws.Range("H2").Value = 0
